$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": insert new row 47 (2020/5/29, serial 43975) before the
# trailing note row, copying the number/date formatting from the row
# directly above (row 46) and filling in the day's figures.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()
$wsAll.Rows(47).Insert(-4121, 0)
$wsAll.Range("A47").Value = 43975
$wsAll.Range("B47").Value = 285
$wsAll.Range("C47").Value = 282
$wsAll.Range("D47").Value = 32
$wsAll.Range("E47").Value = 28
$wsAll.Range("F47").Value = 4
$wsAll.Range("G47").Value = 12
$wsAll.Range("H47").Value = 238
$wsAll.Range("H47").Select()

# ---------------------------------------------------------------------
# Sheet "kobe": insert new row 102 (same date) before the trailing note
# row, copying formatting from row 101, and fill the day's figures.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()
$wsKobe.Rows(102).Insert(-4121, 0)
$wsKobe.Range("A102").Value = 43975
$wsKobe.Range("B102").Value = 0
$wsKobe.Range("C102").Value = 3021
$wsKobe.Range("D102").Value = 0
$wsKobe.Range("E102").Value = 285
$wsKobe.Range("F102").Value = 27
$wsKobe.Range("G102").Value = 24
$wsKobe.Range("H102").Value = 3
$wsKobe.Range("I102").Value = 12
$wsKobe.Range("J102").Value = 229
$wsKobe.Range("G103").Select()

# ---------------------------------------------------------------------
# Sheet "other": insert new row 77 (same date) before the trailing note
# row, copying formatting from row 76, and fill the day's figures.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()
$wsOther.Rows(77).Insert(-4121, 0)
$wsOther.Range("A77").Value = 43975
$wsOther.Range("B77").Value = 0
$wsOther.Range("C77").Value = 14
$wsOther.Range("D77").Value = 5
$wsOther.Range("E77").Value = 4
$wsOther.Range("F77").Value = 1
$wsOther.Range("G77").Value = 0
$wsOther.Range("H77").Value = 9
$wsOther.Range("H77").Select()

# Leave "all" as the active sheet/tab, matching the workbook's original
# tabSelected state.
$wsAll.Activate()
